$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded; insert it as row 12,
# pushing the existing rows (12..135) down to (13..136).
$ws.Rows(12).Insert()

$ws.Range("A12").Value = 5
$ws.Range("B12").Value = 'Macroferia Regional de Talca'
$ws.Range("C12").Value = 'Maule'
$ws.Range("D12").Value = 44503
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 100112017
$ws.Range("G12").Value = 'Apio'
$ws.Range("H12").Value = 'Americana (o)'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 7500
$ws.Range("L12").Value = 7500
$ws.Range("M12").Value = 7500
$ws.Range("N12").Value = '$/docena de matas'
$ws.Range("O12").Value = 'Provincia del Elquí'
$ws.Range("P12").Value = 1250
$ws.Range("Q12").Value = 6
$ws.Range("R12").Value = 'Hortaliza'
